# M1 compatibility resolved & updated the doc
# Update the remaining EOSIO-specific terms in the "Solana" column (B) of the
# Solana vs EVM comparison table to their Solana equivalents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Types of accounts" row: "1. User account / 2. Contract account" -> "1. Account / 2. Program"
$ws.Range("B4").Value = "1. Account" + [char]10 + "2. Program"

# "chain data storage folder" row: "nodeos" -> "test-ledger/"
$ws.Range("B5").Value = "test-ledger/"

# "32 bit" row: "uint32_t" -> "u32"
$ws.Range("B11").Value = "u32"

# "token standard" row: "EOSIO token" -> "Token program"
$ws.Range("B17").Value = "Token program"

# Cells that were actually corrected (no longer red/flagged) now use the
# plain "Segoe UI" style instead of the red "needs updating" style.
# (B6 already carries that "corrected" style in this sheet - copy its format.)
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B11").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# Update viewport / selection to match the author's final cursor position.
$ws.Activate() | Out-Null
$ws.Range("B18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
